$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.409.27'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.66%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.641.51'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -1.47%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +3.86%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.03'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -3.56%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -2.39%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0610'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -1.68%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0891'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +1.36%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.874.67'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.630.93'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -3.79%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +1.84%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -2.17%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.39'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -3.08%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '27.386.16'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.69%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '228.83'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -5.24%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.57'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -1.15%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -0.13%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -3.79%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.57'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +2.45%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.02'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '147.21'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -3.08%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +1.30%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -5.03%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -1.98%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.30%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.415.08'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -3.16%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.64%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.883'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -4.46%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -4.16%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -0.92%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.832'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +5.61%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'FraxShare'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.50'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +1.58%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = 'Aave'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '64.57'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -7.24%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = 'RocketPoolETH'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.784.01'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -1.38%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.67'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -4.28%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '88.00'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -1.18%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = 'BabyDogeCoin'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0₆0107'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0990'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -3.39%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -1.97%  '
$c.NumberFormat = "General"
$c.Style = "Normal"

